$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Extend formatting down to the new row 16 using row 15's format as a template ---
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# --- Shift rows 12-15 down to 13-16 (bottom-up so we don't clobber data before reading it) ---
$ws.Range("A16").Value = $ws.Range("A15").Value()
$ws.Range("B16").Value = $ws.Range("B15").Value()
$ws.Range("A15").Value = $ws.Range("A14").Value()
$ws.Range("B15").Value = $ws.Range("B14").Value()
$ws.Range("A14").Value = $ws.Range("A13").Value()
$ws.Range("B14").Value = $ws.Range("B13").Value()
$ws.Range("A13").Value = $ws.Range("A12").Value()
$ws.Range("B13").Value = $ws.Range("B12").Value()

# --- New "Jurisdiction" row inserted at row 12 ---
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# --- Field value updates ---
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"
